{"js": "async (context) => {\n  // 1) Fix the \"Adress Email:\" typo -> \"Adres Email:\" (keep the existing bold\n  //    formatting that both halves of the original, split run already shared).\n  const emailLabel = context.document.body.search(\"Adress Email:\", { matchCase: true });\n  emailLabel.load(\"text\");\n  await context.sync();\n\n  if (emailLabel.items.length > 0) {\n    const fixedRange = emailLabel.items[0];\n    fixedRange.insertText(\"Adres Email:\", Word.InsertLocation.replace);\n    await context.sync();\n\n    // Re-apply bold explicitly so the corrected (now single) run keeps the\n    // same look the mis-spelled text had.\n    const reloaded = context.document.body.search(\"Adres Email:\", { matchCase: true });\n    reloaded.load(\"text\");\n    await context.sync();\n    if (reloaded.items.length > 0) {\n      reloaded.items[0].font.bold = true;\n      await context.sync();\n    }\n  }\n\n  // 2) Tidy up a sentence that was left split across two runs - re-inserting\n  //    the same text merges it back into a single run.\n  const krok2Sentence =\n    \"przyciskiem \\u201eNowy\\u201d (rys.9).\";\n  const krok2Results = context.document.body.search(krok2Sentence, { matchCase: true });\n  krok2Results.load(\"text\");\n  await context.sync();\n  if (krok2Results.items.length > 0) {\n    krok2Results.items[0].insertText(krok2Sentence, Word.InsertLocation.replace);\n    await context.sync();\n  }\n\n  // 3) Same clean-up for the \"W przypadku udanego dodania nowego cennika...\"\n  //    paragraph, which was split across five separate runs.\n  const successSentence =\n    \"W przypadku udanego dodania nowego cennika, strona zostanie od\u015bwie\u017cona, \" +\n    \"a w formularzu nie zostanie wy\u015bwietlony komunikat i b\u0119dzie mo\u017cna wykonywa\u0107 inne czynno\u015bci.\";\n  const successResults = context.document.body.search(successSentence, { matchCase: true });\n  successResults.load(\"text\");\n  await context.sync();\n  if (successResults.items.length > 0) {\n    successResults.items[0].insertText(successSentence, Word.InsertLocation.replace);\n    await context.sync();\n  }\n};\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Fix the \"Adress Email:\" typo -> \"Adres Email:\" and keep it bold, the\n#    way the original (mis-spelled) text already was.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"Adress Email:\"\n$find.MatchCase = $true\n$find.Execute() | Out-Null\nif ($find.Found) {\n    $rng.Text = \"Adres Email:\"\n    $rng.Font.Bold = 1\n}\n\n# 2) Merge the sentence that was left split across two runs back together.\n#    Word's Range.Text setter is a no-op if the replacement text is\n#    byte-for-byte identical to what is already there, so first swap in a\n#    (trivially different) placeholder and then set the real text - this\n#    forces a real content rewrite and the run gets re-created as one piece.\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.ClearFormatting()\n$find2.Text = \"przyciskiem \u201eNowy\u201d (rys.9).\"\n$find2.MatchCase = $true\n$find2.Execute() | Out-Null\nif ($find2.Found) {\n    $rng2.Text = \"przyciskiem \u201eNowy\u201d (rys.9). \"\n    $rng2.Text = \"przyciskiem \u201eNowy\u201d (rys.9).\"\n}\n\n# 3) Same clean-up for the \"udanego dodania nowego cennika\" paragraph, which\n#    was split across five separate runs.\n$rng3 = $d.Content\n$find3 = $rng3.Find\n$find3.ClearFormatting()\n$find3.Text = \"W przypadku udanego dodania nowego cennika, strona zostanie od\u015bwie\u017cona, a w formularzu nie zostanie wy\u015bwietlony komunikat i b\u0119dzie mo\u017cna wykonywa\u0107 inne czynno\u015bci.\"\n$find3.MatchCase = $true\n$find3.Execute() | Out-Null\nif ($find3.Found) {\n    $rng3.Text = \"W przypadku udanego dodania nowego cennika, strona zostanie od\u015bwie\u017cona, a w formularzu nie zostanie wy\u015bwietlony komunikat i b\u0119dzie mo\u017cna wykonywa\u0107 inne czynno\u015bci. \"\n    $rng3.Text = \"W przypadku udanego dodania nowego cennika, strona zostanie od\u015bwie\u017cona, a w formularzu nie zostanie wy\u015bwietlony komunikat i b\u0119dzie mo\u017cna wykonywa\u0107 inne czynno\u015bci.\"\n}\n"}
